# LeadTemplate.xlsx update
# - Rename hidden lookup sheet "LeadSource" -> "Data"
# - Add Districts / Localities / Pincodes / States / SubDistricts lookup
#   columns (with their own Tables) to the "Data" sheet
# - Point the existing LeadSource list-validation at the renamed sheet and
#   add new list validations for Pincode / Locality / SubDistrict on Sheet1
# - Leave the view/selection state the way the author left it

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$wsData = $wb.Worksheets.Item("LeadSource")

# --- 1. rename the hidden lookup sheet ------------------------------------
$wsData.Name = "Data"

# --- 2. populate the new lookup columns -----------------------------------
# (typed in the same order the author originally entered them so the shared
#  string table / sheet layout comes out the same)

# Districts (column D)
$wsData.Range("D1").Value = "Districts"
$wsData.Range("D2").Value = "Cuddalore"
$wsData.Range("D3").Value = "Vellore"

# Localities data first, header added later (column F)
$wsData.Range("F2").Value = "Tiruvamur"
$wsData.Range("F3").Value = "Kattiyampalayam"

# States data first, header added later (column J)
$wsData.Range("J2").Value = "Tamilnadu"
$wsData.Range("J3").Value = "Kerala"

# headers for States / Pincodes / Localities
$wsData.Range("J1").Value = "States"
$wsData.Range("H1").Value = "Pincodes"
$wsData.Range("F1").Value = "Localities"

# Pincodes data (column H, numeric)
$wsData.Range("H2").Value = 607106
$wsData.Range("H3").Value = 607108

# SubDistricts (column L)
$wsData.Range("L1").Value = "SubDistricts"
$wsData.Range("L2").Value = "Panruti"
$wsData.Range("L3").Value = "Virudhachalam"

# --- 3. turn each new lookup column into its own Table --------------------
$tblDistricts = $wsData.ListObjects.Add(1, $wsData.Range("D1:D3"), $null, 1)
$tblDistricts.Name = "Table3"
$tblDistricts.TableStyle = "TableStyleLight8"

$tblLocalities = $wsData.ListObjects.Add(1, $wsData.Range("F1:F3"), $null, 1)
$tblLocalities.Name = "Table4"
$tblLocalities.TableStyle = "TableStyleLight8"

$tblPincodes = $wsData.ListObjects.Add(1, $wsData.Range("H1:H3"), $null, 1)
$tblPincodes.Name = "Table5"
$tblPincodes.TableStyle = "TableStyleLight8"

$tblStates = $wsData.ListObjects.Add(1, $wsData.Range("J1:J3"), $null, 1)
$tblStates.Name = "Table6"
$tblStates.TableStyle = "TableStyleLight8"

$tblSubDistricts = $wsData.ListObjects.Add(1, $wsData.Range("L1:L3"), $null, 1)
$tblSubDistricts.Name = "Table7"
$tblSubDistricts.TableStyle = "TableStyleLight8"

# --- 4. data validations on Sheet1 -----------------------------------------
# Existing LeadSource dropdown now points at the renamed "Data" sheet.
$ws1.Range("M2").Validation.Delete()
$ws1.Range("M2").Validation.Add(3, $null, 1, "Data!`$A`$2:`$A`$4")

# New Pincode dropdown (no blanks allowed)
$ws1.Range("H2").Validation.Add(3, $null, 1, "Data!`$H`$2:`$H`$3")
$ws1.Range("H2").Validation.IgnoreBlank = $false

# New Locality dropdown
$ws1.Range("I2").Validation.Add(3, $null, 1, "Data!`$F`$2:`$F`$3")

# New SubDistrict dropdown
$ws1.Range("J2").Validation.Add(3, $null, 1, "Data!`$L`$2:`$L`$3")

# --- 5. leave the UI selection state the way the author saved it ----------
$wsData.Activate()
$wsData.Range("K3").Select()

$ws1.Activate()
$ws1.Range("G2").Select()
